# artifacts.xlsx - Items sheet update
# - adjusts stat-mod values on rows 2-3
# - reworks row 4 (Emerald Laced Bow) stat mods + id
# - adds a brand-new artifact row 5 (Ancestral Soldiers Statue / Ice Strength)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Items")

# ---------------------------------------------------------------------------
# Row 2 (Ancestral Finger Bone of The Magi Troth)
# ---------------------------------------------------------------------------
$ws.Range("Q2").Value = 0.3
$ws.Range("R2").ClearContents()
$ws.Range("U2").Value = 0.45
$ws.Range("X2").Value = 0.6
$ws.Range("Z2").Value = 0.6
$ws.Range("AS2").ClearContents()
$ws.Range("AX2").ClearContents()
$ws.Range("AY2").ClearContents()
$ws.Range("BA2").ClearContents()

# ---------------------------------------------------------------------------
# Row 3 (Ancestral Witches Ice Bracelet)
# ---------------------------------------------------------------------------
$ws.Range("Q3").Value = 0.3
$ws.Range("R3").Value = 0.45
$ws.Range("S3").Value = 0.25
$ws.Range("U3").Value = 0.45
$ws.Range("W3").Value = 0.6
$ws.Range("Z3").Value = 0.6
$ws.Range("AS3").ClearContents()
$ws.Range("AX3").ClearContents()
$ws.Range("AY3").ClearContents()
$ws.Range("BA3").ClearContents()

# ---------------------------------------------------------------------------
# Row 4 (Emerald Laced Bow)
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = 1170414
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("Q4").Value = 0.2
$ws.Range("R4").Value = 0.15
$ws.Range("S4").Value = 0.1
$ws.Range("U4").Value = 0.5
$ws.Range("V4").Value = 0.8
$ws.Range("Y4").Value = 0.6
$ws.Range("Z4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 0
$ws.Range("AK4").Value = 0
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = 0
$ws.Range("BB4").Value = 0
$ws.Range("BC4").Value = 0
$ws.Range("BD4").Value = 0
$ws.Range("BE4").Value = 0
$ws.Range("BF4").Value = 0
$ws.Range("BG4").Value = 0

# ---------------------------------------------------------------------------
# Row 5 (new) - Ancestral Soldiers Statue / Ice Strength
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = 2394030
$ws.Range("C5").Value = "Ancestral Soldiers Statue"
$ws.Range("D5").Value = "artifact"
$ws.Range("G5").Value = "A simple soldiers statue from a frozen waste land of depsiar and loss. What secrets could this item hold, what lies will it whisper in the icey cold winds?"
$ws.Range("Q5").Value = 0.25
$ws.Range("S5").Value = 0.3
$ws.Range("T5").Value = 0.4
$ws.Range("U5").Value = 0.15
$ws.Range("V5").Value = 0.1
$ws.Range("AC5").Value = 1
$ws.Range("AV5").Value = 0
$ws.Range("BM5").Value = 0
$ws.Range("BN5").Value = 0
$ws.Range("BO5").Value = 0
$ws.Range("BP5").Value = 0
$ws.Range("BQ5").Value = 0
$ws.Range("BT5").Value = "Ice Strength"
